$d = $word.ActiveDocument

# Remove the stray "_GoBack" bookmark that wrapped the empty tail of the last bullet
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the last bullet paragraph ("Agregar ") - it is the 7th paragraph
$p7 = $d.Paragraphs.Item(7)
$r = $p7.Range

# Range covering the paragraph's content only (excludes the trailing paragraph mark)
$rr = $d.Range($r.Start, $r.End - 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Agregar </w:t></w:r><w:r><w:t>formulario producto</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Poner el logo como componente</w:t></w:r><w:r><w:t>6</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Crear componente botón y agregarlo a los </w:t></w:r><w:r><w:t>distintos formularios</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Elección paleta de colores para el fondo del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Header</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:  #</w:t></w:r><w:r><w:t>f9fcfc</w:t></w:r><w:r><w:t xml:space="preserve">, para el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Footer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: #</w:t></w:r><w:r><w:t>bce5ae</w:t></w:r><w:r><w:t xml:space="preserve">, el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Main</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tendrá un fondo de pantalla de “suculenta”.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Se añaden iconos de Facebook, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>twitter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>whatsapp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rr.InsertXML($xml)
